$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.832.82'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('D3').Value = '3.276.67'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.75'
$ws.Range('E5').Value = '  -1.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.05'
$ws.Range('E6').Value = '  -5.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.580'
$ws.Range('E8').Value = '  +1.52%  '
$ws.Range('D9').Value = '3.268.05'
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('E10').Value = '  -3.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.569'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.15'
$ws.Range('E12').Value = '  -3.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000266'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '683.07'
$ws.Range('E14').Value = '  +7.44%  '
$ws.Range('D15').Value = '3.805.70'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.23'
$ws.Range('E16').Value = '  -2.67%  '
$ws.Range('D17').Value = '66.983.95'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('D19').Value = '3.282.52'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.21'
$ws.Range('E20').Value = '  -3.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.67'
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.881'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.85'
$ws.Range('E23').Value = '  -4.76%  '
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.42'
$ws.Range('E25').Value = '  -2.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.84'
$ws.Range('E26').Value = '  -3.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.67'
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.19'
$ws.Range('E28').Value = '  -2.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.02'
$ws.Range('E29').Value = '  +6.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.31'
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.66'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '572.26'
$ws.Range('E32').Value = '  -3.27%  '
$ws.Range('D33').Value = '3.869.04'
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.75'
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.102'
$ws.Range('E35').Value = '  -2.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.37'
$ws.Range('E37').Value = '  -0.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.27'
$ws.Range('E38').Value = '  -13.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.127'
$ws.Range('E39').Value = '  +0.53%  '
$ws.Range('B40').Value = 'ApeXProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.36'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.57'
$ws.Range('E41').Value = '  -1.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '31.51'
$ws.Range('E42').Value = '  -2.95%  '
$ws.Range('D43').Value = '0.0₃0663'
$ws.Range('E43').Value = '  -5.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.98'
$ws.Range('E44').Value = '  -5.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.325'
$ws.Range('E45').Value = '  -2.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0402'
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.01'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.126'
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('E49').Value = '  +6.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.50'
$ws.Range('E50').Value = '  -1.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '130.06'
$ws.Range('E51').Value = '  -0.54%  '
